$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 16.98068156391676
$ws.Cells.Item(2, 3).Value = 13.36903856161564
$ws.Cells.Item(2, 4).Value = 14.2089050573277
$ws.Cells.Item(2, 5).Value = 14.98716360428322
$ws.Cells.Item(2, 7).Value = 48.75384679447161
$ws.Cells.Item(2, 8).Value = 18.95781677691406
$ws.Cells.Item(2, 9).Value = 26.39963883892306
$ws.Cells.Item(2, 10).Value = 8.869994500130813
$ws.Cells.Item(2, 13).Value = 20.42518549093617
$ws.Cells.Item(2, 14).Value = 19.19047448283403
$ws.Cells.Item(3, 2).Value = 16.55480072689673
$ws.Cells.Item(3, 3).Value = 12.97481776453427
$ws.Cells.Item(3, 4).Value = 14.19956032489629
$ws.Cells.Item(3, 5).Value = 15.00660191405924
$ws.Cells.Item(3, 7).Value = 48.51016893378394
$ws.Cells.Item(3, 8).Value = 18.9783591754774
$ws.Cells.Item(3, 9).Value = 26.45095992346798
$ws.Cells.Item(3, 10).Value = 8.889068909035265
$ws.Cells.Item(3, 13).Value = 20.30350687013461
$ws.Cells.Item(3, 14).Value = 19.26399298600568
$ws.Cells.Item(4, 2).Value = 16.29207405841201
$ws.Cells.Item(4, 3).Value = 12.73041093641418
$ws.Cells.Item(4, 4).Value = 14.19674896343083
$ws.Cells.Item(4, 5).Value = 15.02111888294094
$ws.Cells.Item(4, 7).Value = 48.37644988977072
$ws.Cells.Item(4, 8).Value = 18.99526670356897
$ws.Cells.Item(4, 9).Value = 26.48889933771806
$ws.Cells.Item(4, 10).Value = 8.901495643222651
$ws.Cells.Item(4, 13).Value = 20.23302919380935
$ws.Cells.Item(4, 14).Value = 19.31108547821958
$ws.Cells.Item(5, 2).Value = 16.18487196569146
$ws.Cells.Item(5, 3).Value = 12.63038978752224
$ws.Cells.Item(5, 4).Value = 14.19633992895493
$ws.Cells.Item(5, 5).Value = 15.02768314665851
$ws.Cells.Item(5, 7).Value = 48.3259914439362
$ws.Cells.Item(5, 8).Value = 19.00323353256422
$ws.Cells.Item(5, 9).Value = 26.50596955094305
$ws.Cells.Item(5, 10).Value = 8.906739868733629
$ws.Cells.Item(5, 13).Value = 20.2053965539254
$ws.Cells.Item(5, 14).Value = 19.33076857917101
$ws.Cells.Item(6, 2).Value = 16.1670678924775
$ws.Cells.Item(6, 3).Value = 12.61376088897986
$ws.Cells.Item(6, 4).Value = 14.19631651278078
$ws.Cells.Item(6, 5).Value = 15.02881228061327
$ws.Cells.Item(6, 7).Value = 48.31785733244855
$ws.Cells.Item(6, 8).Value = 19.00462136403905
$ws.Cells.Item(6, 9).Value = 26.50890105537463
$ws.Cells.Item(6, 10).Value = 8.907621566991594
$ws.Cells.Item(6, 13).Value = 20.20087450215101
$ws.Cells.Item(6, 14).Value = 19.33406674114801
$ws.Cells.Item(7, 2).Value = 16.29062862432858
$ws.Cells.Item(7, 3).Value = 12.72906350087582
$ws.Cells.Item(7, 4).Value = 14.19674046382738
$ws.Cells.Item(7, 5).Value = 15.02120478628404
$ws.Cells.Item(7, 7).Value = 48.37575301873541
$ws.Cells.Item(7, 8).Value = 18.99536979111868
$ws.Cells.Item(7, 9).Value = 26.48912304481613
$ws.Cells.Item(7, 10).Value = 8.901565638349554
$ws.Cells.Item(7, 13).Value = 20.23265209765535
$ws.Cells.Item(7, 14).Value = 19.31134893494222
$ws.Cells.Item(8, 2).Value = 16.83419475278802
$ws.Cells.Item(8, 3).Value = 13.23369778501249
$ws.Cells.Item(8, 4).Value = 14.20507621623986
$ws.Cells.Item(8, 5).Value = 14.99332968567645
$ws.Cells.Item(8, 7).Value = 48.66654966773219
$ws.Cells.Item(8, 8).Value = 18.96400680154187
$ws.Cells.Item(8, 9).Value = 26.41599624228685
$ws.Cells.Item(8, 10).Value = 8.876423232157551
$ws.Cells.Item(8, 13).Value = 20.38236637530688
$ws.Cells.Item(8, 14).Value = 19.21541960319883
$ws.Cells.Item(9, 2).Value = 17.88317145974918
$ws.Cells.Item(9, 3).Value = 14.1977216273788
$ws.Cells.Item(9, 4).Value = 14.24459131969207
$ws.Cells.Item(9, 5).Value = 14.95918072875848
$ws.Cells.Item(9, 7).Value = 49.36109862370121
$ws.Cells.Item(9, 8).Value = 18.93669737478534
$ws.Cells.Item(9, 9).Value = 26.3238854971813
$ws.Cells.Item(9, 10).Value = 8.832771911403153
$ws.Cells.Item(9, 13).Value = 20.70845732325498
$ws.Cells.Item(9, 14).Value = 19.04270909964912
$ws.Cells.Item(10, 2).Value = 18.63453197990462
$ws.Cells.Item(10, 3).Value = 14.88196137064332
$ws.Cells.Item(10, 4).Value = 14.28765491954683
$ws.Cells.Item(10, 5).Value = 14.94663150247069
$ws.Cells.Item(10, 7).Value = 49.94421215007688
$ws.Cells.Item(10, 8).Value = 18.93761925235812
$ws.Cells.Item(10, 9).Value = 26.2878467398828
$ws.Cells.Item(10, 10).Value = 8.804119598166247
$ws.Cells.Item(10, 13).Value = 20.96629011585312
$ws.Cells.Item(10, 14).Value = 18.9250962090979
$ws.Cells.Item(11, 2).Value = 18.97044005546319
$ws.Cells.Item(11, 3).Value = 15.18646698449157
$ws.Cells.Item(11, 4).Value = 14.31026354187414
$ws.Cells.Item(11, 5).Value = 14.94365004423478
$ws.Cells.Item(11, 7).Value = 50.22454282532425
$ws.Cells.Item(11, 8).Value = 18.9426170733602
$ws.Cells.Item(11, 9).Value = 26.27838949218908
$ws.Cells.Item(11, 10).Value = 8.791821310261595
$ws.Cells.Item(11, 13).Value = 21.0871890588365
$ws.Cells.Item(11, 14).Value = 18.87358094794005
$ws.Cells.Item(12, 2).Value = 19.09666423892909
$ws.Cells.Item(12, 3).Value = 15.3006895745246
$ws.Cells.Item(12, 4).Value = 14.3192556639958
$ws.Cells.Item(12, 5).Value = 14.94291329347112
$ws.Cells.Item(12, 7).Value = 50.33279163045654
$ws.Cells.Item(12, 8).Value = 18.94516893298768
$ws.Cells.Item(12, 9).Value = 26.27581048640169
$ws.Cells.Item(12, 10).Value = 8.787269632228851
$ws.Cells.Item(12, 13).Value = 21.1334574791892
$ws.Cells.Item(12, 14).Value = 18.85435741461328
$ws.Cells.Item(13, 2).Value = 19.06952500976613
$ws.Cells.Item(13, 3).Value = 15.27613975888995
$ws.Cells.Item(13, 4).Value = 14.31729995139555
$ws.Cells.Item(13, 5).Value = 14.94305452021194
$ws.Cells.Item(13, 7).Value = 50.30938641582546
$ws.Cells.Item(13, 8).Value = 18.94459001249524
$ws.Cells.Item(13, 9).Value = 26.27632127703962
$ws.Cells.Item(13, 10).Value = 8.788245235243581
$ws.Cells.Item(13, 13).Value = 21.12347162992164
$ws.Cells.Item(13, 14).Value = 18.85848493314241
$ws.Cells.Item(14, 2).Value = 18.98084488124889
$ws.Cells.Item(14, 3).Value = 15.19588655221782
$ws.Cells.Item(14, 4).Value = 14.3109947176743
$ws.Cells.Item(14, 5).Value = 14.94358156983128
$ws.Cells.Item(14, 7).Value = 50.23340699669632
$ws.Cells.Item(14, 8).Value = 18.94281379959348
$ws.Cells.Item(14, 9).Value = 26.27815720426834
$ws.Cells.Item(14, 10).Value = 8.791444730476403
$ws.Cells.Item(14, 13).Value = 21.09098602588274
$ws.Cells.Item(14, 14).Value = 18.87199373014071
$ws.Cells.Item(15, 2).Value = 18.92639483929864
$ws.Cells.Item(15, 3).Value = 15.14658434965851
$ws.Cells.Item(15, 4).Value = 14.30718856559833
$ws.Cells.Item(15, 5).Value = 14.94395548705963
$ws.Cells.Item(15, 7).Value = 50.18713780641723
$ws.Cells.Item(15, 8).Value = 18.94181169652837
$ws.Cells.Item(15, 9).Value = 26.27941241544044
$ws.Cells.Item(15, 10).Value = 8.793418230581146
$ws.Cells.Item(15, 13).Value = 21.0711500500526
$ws.Cells.Item(15, 14).Value = 18.88030521901643
$ws.Cells.Item(16, 2).Value = 18.61244973189319
$ws.Cells.Item(16, 3).Value = 14.86191510675684
$ws.Cells.Item(16, 4).Value = 14.28623786091424
$ws.Cells.Item(16, 5).Value = 14.94688123806415
$ws.Cells.Item(16, 7).Value = 49.9261889636503
$ws.Cells.Item(16, 8).Value = 18.93738485231983
$ws.Cells.Item(16, 9).Value = 26.28860479930614
$ws.Cells.Item(16, 10).Value = 8.804938092833572
$ws.Cells.Item(16, 13).Value = 20.95845896605136
$ws.Cells.Item(16, 14).Value = 18.92850270676447
$ws.Cells.Item(17, 2).Value = 18.41825029425922
$ws.Cells.Item(17, 3).Value = 14.68546314770684
$ws.Cells.Item(17, 4).Value = 14.27415621228592
$ws.Cells.Item(17, 5).Value = 14.94937471724276
$ws.Cells.Item(17, 7).Value = 49.76991438124128
$ws.Cells.Item(17, 8).Value = 18.93584257994667
$ws.Cells.Item(17, 9).Value = 26.29602429104001
$ws.Cells.Item(17, 10).Value = 8.812193330726377
$ws.Cells.Item(17, 13).Value = 20.89022832966113
$ws.Cells.Item(17, 14).Value = 18.95857818352165
$ws.Cells.Item(18, 2).Value = 18.30600466643457
$ws.Cells.Item(18, 3).Value = 14.58334266266182
$ws.Cells.Item(18, 4).Value = 14.26749154000436
$ws.Cells.Item(18, 5).Value = 14.95106560666916
$ws.Cells.Item(18, 7).Value = 49.68145250087013
$ws.Cells.Item(18, 8).Value = 18.93538642632702
$ws.Cells.Item(18, 9).Value = 26.30094447959991
$ws.Cells.Item(18, 10).Value = 8.816435630982678
$ws.Cells.Item(18, 13).Value = 20.85132561893828
$ws.Cells.Item(18, 14).Value = 18.97606396476501
$ws.Cells.Item(19, 2).Value = 18.26791036773654
$ws.Cells.Item(19, 3).Value = 14.54866191743064
$ws.Cells.Item(19, 4).Value = 14.26528393280765
$ws.Cells.Item(19, 5).Value = 14.95168219703964
$ws.Cells.Item(19, 7).Value = 49.65174740276315
$ws.Cells.Item(19, 8).Value = 18.93530595589293
$ws.Cells.Item(19, 9).Value = 26.30272231832636
$ws.Cells.Item(19, 10).Value = 8.817883911686216
$ws.Cells.Item(19, 13).Value = 20.83821349840641
$ws.Cells.Item(19, 14).Value = 18.98201654441386
$ws.Cells.Item(20, 2).Value = 18.43898076898371
$ws.Cells.Item(20, 3).Value = 14.70431283158885
$ws.Cells.Item(20, 4).Value = 14.27541291590703
$ws.Cells.Item(20, 5).Value = 14.94908271292235
$ws.Cells.Item(20, 7).Value = 49.7864033048427
$ws.Cells.Item(20, 8).Value = 18.93596215087351
$ws.Cells.Item(20, 9).Value = 26.2951668864606
$ws.Cells.Item(20, 10).Value = 8.811413830979674
$ws.Cells.Item(20, 13).Value = 20.89745646014731
$ws.Cells.Item(20, 14).Value = 18.95535723889071
$ws.Cells.Item(21, 2).Value = 19.00691988195242
$ws.Cells.Item(21, 3).Value = 15.21948922964892
$ws.Cells.Item(21, 4).Value = 14.31283505695617
$ws.Cells.Item(21, 5).Value = 14.94341611676237
$ws.Cells.Item(21, 7).Value = 50.25566778071766
$ws.Cells.Item(21, 8).Value = 18.94331761836773
$ws.Cells.Item(21, 9).Value = 26.27759071239526
$ws.Cells.Item(21, 10).Value = 8.790502103269345
$ws.Cells.Item(21, 13).Value = 21.10051489127909
$ws.Cells.Item(21, 14).Value = 18.86801816929784
$ws.Cells.Item(22, 2).Value = 19.37235377981644
$ws.Cells.Item(22, 3).Value = 15.5498048896854
$ws.Cells.Item(22, 4).Value = 14.33980131497398
$ws.Cells.Item(22, 5).Value = 14.94199900554969
$ws.Cells.Item(22, 7).Value = 50.57452540006636
$ws.Cells.Item(22, 8).Value = 18.95196779936038
$ws.Cells.Item(22, 9).Value = 26.27194732216726
$ws.Cells.Item(22, 10).Value = 8.777449352348729
$ws.Cells.Item(22, 13).Value = 21.23604607880781
$ws.Cells.Item(22, 14).Value = 18.81259270343472
$ws.Cells.Item(23, 2).Value = 19.17788094349501
$ws.Cells.Item(23, 3).Value = 15.37412831967177
$ws.Cells.Item(23, 4).Value = 14.3251805977811
$ws.Cells.Item(23, 5).Value = 14.94254615613212
$ws.Cells.Item(23, 7).Value = 50.40325732769492
$ws.Cells.Item(23, 8).Value = 18.94699922287755
$ws.Cells.Item(23, 9).Value = 26.27442319932586
$ws.Cells.Item(23, 10).Value = 8.784359772591561
$ws.Cells.Item(23, 13).Value = 21.16346335509411
$ws.Cells.Item(23, 14).Value = 18.84202335466615
$ws.Cells.Item(24, 2).Value = 18.42961037330967
$ws.Cells.Item(24, 3).Value = 14.69579298595422
$ws.Cells.Item(24, 4).Value = 14.27484388397438
$ws.Cells.Item(24, 5).Value = 14.94921392646162
$ws.Cells.Item(24, 7).Value = 49.77894435073186
$ws.Cells.Item(24, 8).Value = 18.93590675186624
$ws.Cells.Item(24, 9).Value = 26.29555248029992
$ws.Cells.Item(24, 10).Value = 8.811766021104585
$ws.Cells.Item(24, 13).Value = 20.8941876107001
$ws.Cells.Item(24, 14).Value = 18.95681282043231
$ws.Cells.Item(25, 2).Value = 17.60217479637503
$ws.Cells.Item(25, 3).Value = 13.9405806652558
$ws.Cells.Item(25, 4).Value = 14.23142877994076
$ws.Cells.Item(25, 5).Value = 14.96621867096778
$ws.Cells.Item(25, 7).Value = 49.16017684280662
$ws.Cells.Item(25, 8).Value = 18.94040842200774
$ws.Cells.Item(25, 9).Value = 26.34327304864178
$ws.Cells.Item(25, 10).Value = 8.843978467210217
$ws.Cells.Item(25, 13).Value = 20.616917103144
$ws.Cells.Item(25, 14).Value = 19.08779437754183
